# Update the "想去人数" (want-to-go count) figures that were refreshed by the
# gh-pages data generation job. The same three rows are present (duplicated)
# on both the "展览" sheet and the "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 49
    $ws.Range("F3").Value = 283
    $ws.Range("F4").Value = 19
}
